# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates to the Tonberry Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 901.05884
$ws.Range("I15").Value = 901.05884
$ws.Range("K15").Value = 2703.17652
$ws.Range("M15").Value = -2534.17652
$ws.Range("H41").Value = 653.2222
$ws.Range("J41").Value = 808.3333
$ws.Range("L41").Value = 808.3333
$ws.Range("N41").Value = -1688.3333
$ws.Range("H98").Value = 1533.1052
$ws.Range("I98").Value = 1595.8235
$ws.Range("K98").Value = 1595.8235
$ws.Range("M98").Value = -97.82349999999997
$ws.Range("H112").Value = 6376.8887
$ws.Range("J112").Value = 6376.8887
$ws.Range("L112").Value = 19130.6661
$ws.Range("N112").Value = -21346.6661
$ws.Range("H122").Value = 1533.1052
$ws.Range("I122").Value = 1595.8235
$ws.Range("K122").Value = 4787.470499999999
$ws.Range("M122").Value = -2337.470499999999
$ws.Range("H132").Value = 1465.7931
$ws.Range("I132").Value = 1465.7931
$ws.Range("K132").Value = 4397.379300000001
$ws.Range("M132").Value = -1867.379300000001
$ws.Range("H137").Value = 3455.5
$ws.Range("I137").Value = 1455.7142
$ws.Range("J137").Value = 4728.091
$ws.Range("K137").Value = 4367.142599999999
$ws.Range("L137").Value = 14184.273
$ws.Range("M137").Value = -1817.142599999999
$ws.Range("N137").Value = -19284.273
$ws.Range("H138").Value = 2404.6584
$ws.Range("I138").Value = 2213.8
$ws.Range("J138").Value = 2702.875
$ws.Range("K138").Value = 6641.400000000001
$ws.Range("L138").Value = 8108.625
$ws.Range("M138").Value = -1501.400000000001
$ws.Range("N138").Value = -18388.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3229.9565
$ws.Range("I32").Value = 1742.4728
$ws.Range("J32").Value = 9073.643
$ws.Range("K32").Value = 1742.4728
$ws.Range("L32").Value = 9073.643
$ws.Range("M32").Value = -1455.4728
$ws.Range("N32").Value = -9647.643
$ws.Range("H88").Value = 23383.8
$ws.Range("I88").Value = 2626.3333
$ws.Range("J88").Value = 32279.857
$ws.Range("K88").Value = 2626.3333
$ws.Range("L88").Value = 32279.857
$ws.Range("M88").Value = -2220.3333
$ws.Range("N88").Value = -33091.857
$ws.Range("H91").Value = 23383.8
$ws.Range("I91").Value = 2626.3333
$ws.Range("J91").Value = 32279.857
$ws.Range("K91").Value = 2626.3333
$ws.Range("L91").Value = 32279.857
$ws.Range("M91").Value = -1222.3333
$ws.Range("N91").Value = -35087.857
$ws.Range("H132").Value = 1443.7391
$ws.Range("I132").Value = 1105.1428
$ws.Range("K132").Value = 3315.4284
$ws.Range("M132").Value = -785.4284000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10832
$ws.Range("H80").Value = 6879.7334
$ws.Range("J80").Value = 8588.833000000001
$ws.Range("L80").Value = 8588.833000000001
$ws.Range("N80").Value = -10584.833
$ws.Range("H81").Value = 60000
$ws.Range("J81").Value = 60000
$ws.Range("L81").Value = 60000
$ws.Range("N81").Value = -62122
$ws.Range("H83").Value = 6879.7334
$ws.Range("J83").Value = 8588.833000000001
$ws.Range("L83").Value = 42944.165
$ws.Range("N83").Value = -52928.165
$ws.Range("H84").Value = 60000
$ws.Range("J84").Value = 60000
$ws.Range("L84").Value = 180000
$ws.Range("N84").Value = -190608
$ws.Range("H134").Value = 9309.467000000001
$ws.Range("I134").Value = 9103.200000000001
$ws.Range("K134").Value = 27309.6
$ws.Range("M134").Value = -24774.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9576.923000000001
$ws.Range("J4").Value = 9576.923000000001
$ws.Range("L4").Value = 9576.923000000001
$ws.Range("N4").Value = -9800.923000000001
$ws.Range("H7").Value = 105.22222
$ws.Range("I7").Value = 123.666664
$ws.Range("J7").Value = 68.333336
$ws.Range("K7").Value = 123.666664
$ws.Range("L7").Value = 68.333336
$ws.Range("M7").Value = -10.666664
$ws.Range("N7").Value = -294.333336
$ws.Range("H31").Value = 2743.9412
$ws.Range("I31").Value = 1917.4584
$ws.Range("K31").Value = 1917.4584
$ws.Range("M31").Value = -1622.4584
$ws.Range("H34").Value = 2743.9412
$ws.Range("I34").Value = 1917.4584
$ws.Range("K34").Value = 1917.4584
$ws.Range("M34").Value = -1715.4584
$ws.Range("H43").Value = 38999
$ws.Range("J43").Value = 38999
$ws.Range("L43").Value = 38999
$ws.Range("N43").Value = -39367
$ws.Range("H101").Value = 38999
$ws.Range("J101").Value = 38999
$ws.Range("L101").Value = 38999
$ws.Range("N101").Value = -45489
$ws.Range("H107").Value = 379.92307
$ws.Range("I107").Value = 379.92307
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 379.92307
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1540.07693
$ws.Range("H122").Value = 1329.8928
$ws.Range("I122").Value = 1332.1666
$ws.Range("J122").Value = 1325.8
$ws.Range("K122").Value = 3996.4998
$ws.Range("L122").Value = 3977.4
$ws.Range("M122").Value = -1546.4998
$ws.Range("N122").Value = -8877.4
$ws.Range("H134").Value = 928.6429000000001
$ws.Range("I134").Value = 928.6429000000001
$ws.Range("K134").Value = 2785.9287
$ws.Range("M134").Value = -250.9287000000004
$ws.Range("N107").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H99").Value = 389.33334
$ws.Range("I99").Value = 389.33334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1168.00002
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 1077.99998
$ws.Range("H122").Value = 740.0769
$ws.Range("J122").Value = 976.125
$ws.Range("L122").Value = 8785.125
$ws.Range("N122").Value = -13685.125
$ws.Range("H131").Value = 779.78
$ws.Range("I131").Value = 537
$ws.Range("J131").Value = 795.2766
$ws.Range("K131").Value = 1611
$ws.Range("L131").Value = 2385.8298
$ws.Range("M131").Value = 3429
$ws.Range("N131").Value = -12465.8298
$ws.Range("M63").Value = $null
$ws.Range("M66").Value = $null
$ws.Range("N99").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7500
$ws.Range("I70").Value = 7500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7230
$ws.Range("H73").Value = 7500
$ws.Range("I73").Value = 7500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6564
$ws.Range("H102").Value = 2556.842
$ws.Range("I102").Value = 2560.5454
$ws.Range("J102").Value = 2551.75
$ws.Range("K102").Value = 2560.5454
$ws.Range("L102").Value = 2551.75
$ws.Range("M102").Value = -938.5454
$ws.Range("N102").Value = -5795.75
$ws.Range("H126").Value = 2573151
$ws.Range("J126").Value = 113257
$ws.Range("L126").Value = 339771
$ws.Range("N126").Value = -344711
$ws.Range("H132").Value = 1284892.6
$ws.Range("I132").Value = 1833173.2
$ws.Range("J132").Value = 5571.3335
$ws.Range("K132").Value = 5499519.6
$ws.Range("L132").Value = 16714.0005
$ws.Range("M132").Value = -5496989.6
$ws.Range("N132").Value = -21774.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6166.222
$ws.Range("I7").Value = 1474.75
$ws.Range("K7").Value = 1474.75
$ws.Range("M7").Value = -1362.75
$ws.Range("H40").Value = 10666.5
$ws.Range("H61").Value = 3899.8
$ws.Range("I61").Value = 3166.3333
$ws.Range("K61").Value = 3166.3333
$ws.Range("M61").Value = -2964.3333
$ws.Range("H113").Value = 3899.8
$ws.Range("I113").Value = 3166.3333
$ws.Range("K113").Value = 3166.3333
$ws.Range("M113").Value = -996.3332999999998
$ws.Range("H122").Value = 5570.2
$ws.Range("I122").Value = 3994.125
$ws.Range("J122").Value = 7371.4287
$ws.Range("K122").Value = 11982.375
$ws.Range("L122").Value = 22114.2861
$ws.Range("M122").Value = -9532.375
$ws.Range("N122").Value = -27014.2861
$ws.Range("H126").Value = 6166.222
$ws.Range("I126").Value = 1474.75
$ws.Range("K126").Value = 4424.25
$ws.Range("M126").Value = -1954.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29703.334
$ws.Range("J70").Value = 29703.334
$ws.Range("L70").Value = 29703.334
$ws.Range("N70").Value = -30333.334
$ws.Range("H73").Value = 29703.334
$ws.Range("J73").Value = 29703.334
$ws.Range("L73").Value = 29703.334
$ws.Range("N73").Value = -31887.334
$ws.Range("H81").Value = 1998.75
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939
$ws.Range("H84").Value = 1998.75
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696
$ws.Range("H100").Value = 722
$ws.Range("I100").Value = 333
$ws.Range("K100").Value = 666
$ws.Range("M100").Value = -125
$ws.Range("H107").Value = 767
$ws.Range("I107").Value = 767
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2301
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -381
$ws.Range("H126").Value = 7158.636
$ws.Range("J126").Value = 7811
$ws.Range("L126").Value = 23433
$ws.Range("N126").Value = -28373
$ws.Range("H132").Value = 1920.0476
$ws.Range("I132").Value = 1651.3334
$ws.Range("K132").Value = 4954.0002
$ws.Range("M132").Value = -2424.0002
$ws.Range("N107").Value = $null
